# Natmi LR-pairs (Bmp2-Bmpr2) update following Dr Hou advice.
#
# The sending/target cluster set grows from {FAPs, sCs} to {FAPs, ECs, sCs},
# so the 2x2 (FAPs/sCs) matrix of result rows becomes a full 3x3
# (FAPs/ECs/sCs) matrix, and every Ligand/Receptor expression statistic is
# recomputed for the new cluster composition.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusters = @("FAPs", "ECs", "sCs")

# Per sending-cluster ligand (Bmp2) stats: E,F,G,H,I,J
$ligandStats = @{
    "FAPs" = @(2, 0.6666666666666666, 1.445484, 4.336452, 0.1286708197254238, 0.1286708197254238)
    "ECs"  = @(3, 1,                  6.292313, 18.876939, 0.5601148623429528, 0.5601148623429528)
    "sCs"  = @(3, 1,                  3.496172, 10.488516, 0.3112143179316233, 0.3112143179316232)
}

# Per target-cluster receptor (Bmpr2) stats: K,L,M,N,O,P
$receptorStats = @{
    "FAPs" = @(3, 1, 40.70766766666667, 122.123003, 0.3776398983502007, 0.3776398983502007)
    "ECs"  = @(3, 1, 39.715023,         119.145069, 0.3684312589831062, 0.3684312589831062)
    "sCs"  = @(3, 1, 27.37224266666666, 82.11672799999999, 0.253928842666693, 0.253928842666693)
}

# Per (sending, target) edge stats: Q,R,S,T
$edgeStats = @{
    "FAPs|FAPs" = @(58.84228228948399,  529.580540605356,   0.04859123528174605, 0.04859123528174605)
    "FAPs|ECs"  = @(57.40743030613199,  516.666872755188,   0.04740635210582619, 0.04740635210582619)
    "FAPs|sCs"  = @(39.56613881878399,  356.0952493690559,  0.03267323233785157, 0.03267323233785157)
    "ECs|FAPs"  = @(256.1453864586463,  2305.308478127817,  0.2115217196796294,  0.2115217196796293)
    "ECs|ECs"   = @(249.899355518199,   2249.094199663791,  0.2063638239081633,  0.2063638239081633)
    "ECs|sCs"   = @(172.2347183706213,  1550.112465335592,  0.1422293187551601,  0.1422293187551601)
    "sCs|FAPs"  = @(142.3210078815053,  1280.889070933548,  0.1175269433888253,  0.1175269433888252)
    "sCs|ECs"   = @(138.850551391956,   1249.654962527604,  0.1146610829691166,  0.1146610829691166)
    "sCs|sCs"   = @(95.69806838840533,  861.282615495648,   0.07902629157368135, 0.07902629157368134)
}

$row = 2
foreach ($sending in $clusters) {
    foreach ($target in $clusters) {
        $ws.Cells.Item($row, 1).Value = $sending
        $ws.Cells.Item($row, 2).Value = "Bmp2"
        $ws.Cells.Item($row, 3).Value = "Bmpr2"
        $ws.Cells.Item($row, 4).Value = $target

        $lig = $ligandStats[$sending]
        $ws.Cells.Item($row, 5).Value  = $lig[0]
        $ws.Cells.Item($row, 6).Value  = $lig[1]
        $ws.Cells.Item($row, 7).Value  = $lig[2]
        $ws.Cells.Item($row, 8).Value  = $lig[3]
        $ws.Cells.Item($row, 9).Value  = $lig[4]
        $ws.Cells.Item($row, 10).Value = $lig[5]

        $rec = $receptorStats[$target]
        $ws.Cells.Item($row, 11).Value = $rec[0]
        $ws.Cells.Item($row, 12).Value = $rec[1]
        $ws.Cells.Item($row, 13).Value = $rec[2]
        $ws.Cells.Item($row, 14).Value = $rec[3]
        $ws.Cells.Item($row, 15).Value = $rec[4]
        $ws.Cells.Item($row, 16).Value = $rec[5]

        $edge = $edgeStats["$sending|$target"]
        $ws.Cells.Item($row, 17).Value = $edge[0]
        $ws.Cells.Item($row, 18).Value = $edge[1]
        $ws.Cells.Item($row, 19).Value = $edge[2]
        $ws.Cells.Item($row, 20).Value = $edge[3]

        $row++
    }
}
